$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.118598818778992
$ws.Range("B1").Value = 2.262534618377686
$ws.Range("C1").Value = 10.79081344604492
$ws.Range("D1").Value = 1.752775907516479
$ws.Range("E1").Value = 1.289687871932983
